# Fix cumulative/total formulas on the "2021" and "2020" report sheets.
# ("2022" sheet already has these SUM formulas from an earlier fix.)
#
# For each affected sheet:
#   - Column N (Cumulative Year), rows 2-6: N<r> = SUM(B<r>:M<r>)
#   - Row 7 (Total), columns B-N: <col>7 = SUM(<col>2:<col>6)
#
# The "2020" sheet additionally gets corrected actuals for April/May
# (columns E/F) on the Premier, Suivi and Commission Taxes rows.

$wb = $excel.ActiveWorkbook

$sheetNames = @("2021", "2020")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    for ($row = 2; $row -le 6; $row++) {
        $ws.Range("N$row").Formula = "=SUM(B$row`:M$row)"
    }

    $cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N")
    foreach ($col in $cols) {
        $ws.Range("$col" + "7").Formula = "=SUM(" + $col + "2:" + $col + "6)"
    }
}

# Corrected historical values on the "2020" sheet
$ws2020 = $wb.Worksheets.Item("2020")
$ws2020.Range("E3").Value = 3491
$ws2020.Range("F3").Value = 3491
$ws2020.Range("E4").Value = 2875
$ws2020.Range("F4").Value = 2875
$ws2020.Range("E6").Value = 953.3085
$ws2020.Range("F6").Value = 953.3085

$wb.Save()
$excel.Quit()
